$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking text (e.g. "1.004", "27.643.64")
# that must remain plain text, matching the source data exactly.
# Force text format, assign the value, then restore default styling
# so no stray number-format style lingers on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.643.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.866.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4706"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3927"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.87"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08034"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.022"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.97"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.71%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.942"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.829.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.130"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.006"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001045"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "86.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06644"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.672.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.493"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.308"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.088.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.088"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.549"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9663"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09498"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.446"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.590"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.32%  "
$ws.Range("E36").Value = "  +0.99%  "
$ws.Range("E37").Value = "  +1.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02252"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.227"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.107"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6033"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.002"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1893"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.258"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5706"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.40%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.935"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.53%  "
$ws.Range("B49").Value = "PancakeSwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.380"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("E50").Value = "  +1.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "114.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.01%  "
